$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.425.57"
$ws.Range("E2").Value = "  -4.54%  "
$ws.Range("D3").Value = "2.927.92"
$ws.Range("E3").Value = "  -1.14%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "552.12"
$ws.Range("E5").Value = "  -3.38%  "
$ws.Range("D6").Value = "130.45"
$ws.Range("E6").Value = "  +6.12%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "0.516"
$ws.Range("E8").Value = "  +3.84%  "
$ws.Range("D9").Value = "2.931.34"
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("E10").Value = "  -2.44%  "
$ws.Range("D11").Value = "4.81"
$ws.Range("E11").Value = "  -4.26%  "
$ws.Range("D12").Value = "0.444"
$ws.Range("E12").Value = "  +2.20%  "
$ws.Range("D13").Value = "0.0000222"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("D14").Value = "32.51"
$ws.Range("E14").Value = "  +0.48%  "
$ws.Range("E15").Value = "  +1.43%  "
$ws.Range("D16").Value = "3.414.93"
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("D17").Value = "6.75"
$ws.Range("E17").Value = "  +9.97%  "
$ws.Range("D18").Value = "2.924.28"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("D19").Value = "57.446.82"
$ws.Range("E19").Value = "  -4.51%  "
$ws.Range("D20").Value = "416.70"
$ws.Range("E20").Value = "  -3.00%  "
$ws.Range("D21").Value = "13.06"
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("D22").Value = "0.678"
$ws.Range("E22").Value = "  +3.38%  "
$ws.Range("D23").Value = "6.93"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").Value = "12.88"
$ws.Range("E24").Value = "  +0.77%  "
$ws.Range("D25").Value = "79.06"
$ws.Range("E25").Value = "  +1.03%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D29").Value = "7.46"
$ws.Range("E29").Value = "  +4.99%  "
$ws.Range("D30").Value = "1.99"
$ws.Range("E30").Value = "  +6.52%  "
$ws.Range("D31").Value = "6.15"
$ws.Range("E31").Value = "  +1.60%  "
$ws.Range("E32").Value = "  +12.52%  "
$ws.Range("D33").Value = "24.99"
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("D34").Value = "5.61"
$ws.Range("E34").Value = "  +1.55%  "
$ws.Range("E35").Value = "  -2.94%  "
$ws.Range("D36").Value = "0.934"
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("D37").Value = "48.65"
$ws.Range("E37").Value = "  -1.60%  "
$ws.Range("D38").Value = "0.0₃0677"
$ws.Range("E38").Value = "  +5.00%  "
$ws.Range("D39").Value = "8.38"
$ws.Range("E39").Value = "  +7.03%  "
$ws.Range("E40").Value = "  +6.88%  "
$ws.Range("E41").Value = "  -1.66%  "
$ws.Range("D42").Value = "0.108"
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("D43").Value = "376.50"
$ws.Range("E43").Value = "  +0.66%  "
$ws.Range("D44").Value = "2.633.09"
$ws.Range("E44").Value = "  +0.81%  "
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").Value = "0.239"
$ws.Range("E46").Value = "  +2.24%  "
$ws.Range("D47").Value = "121.01"
$ws.Range("E47").Value = "  +1.70%  "
$ws.Range("E48").Value = "  +3.38%  "
$ws.Range("E49").Value = "  +1.30%  "
$ws.Range("D50").Value = "23.21"
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("D51").Value = "1.99"
$ws.Range("E51").Value = "  +1.76%  "
